$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the duplicated bold "Play Brazilian Beauty Free..." paragraph
#    that used to sit right before the final italic paragraph.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 2; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Play Brazilian Beauty Free: Captivating Brazilian Slot`r") {
        $p.Range.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2) Replace the text of the final italic paragraph with the new
#    image-generation prompt (keep its italic formatting).
#    We use Find to locate the range, then assign .Text directly so
#    the runtime's smart-quote autocorrect doesn't mangle the straight
#    quotes around "Brazilian Beauty". This runs before the new "Meta
#    description" paragraph is inserted so Find's first (and only) match
#    is this trailing paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Read our review of Brazilian Beauty, an online slot game with a captivating Brazilian theme. Play for free and enjoy bonus features like free spins.") | Out-Null
$rng.Text = "Create a cartoon-style image featuring a happy Maya warrior with glasses for the game ""Brazilian Beauty"". The Maya warrior should be portrayed with bright colors, holding maracas, and with a happy expression. The background should showcase the iconic imagery of Brazil, such as the Cristo Redentor, the beaches, and the rainforest. The image should be eye-catching and colorful, giving players a glimpse of the fun and excitement they can experience while playing this game."

# ---------------------------------------------------------------------------
# 3) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$title = $d.Paragraphs.First
$newPara = $title.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start

$boldRun = $d.Range($metaStart, $metaStart)
$boldRun.InsertAfter("Meta description")

$restRun = $d.Range($boldRun.End, $boldRun.End)
$restRun.InsertAfter(": Read our review of Brazilian Beauty, an online slot game with a captivating Brazilian theme. Play for free and enjoy bonus features like free spins.")

$boldRun.Font.Bold = $true

Write-Output "done"
